$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 815
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 815
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 2445
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -4193

$ws.Range("H72").Value = 815
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 815
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 7335
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -16071

$ws.Range("H98").Value = 1726.5555
$ws.Range("I98").Value = 1409.6923
$ws.Range("J98").Value = 2550.4
$ws.Range("K98").Value = 1409.6923
$ws.Range("L98").Value = 2550.4
$ws.Range("M98").Value = 88.30770000000007
$ws.Range("N98").Value = -5546.4

$ws.Range("H111").Value = 644
$ws.Range("I111").Value = 637.7778
$ws.Range("J111").Value = 700
$ws.Range("K111").Value = 1913.3334
$ws.Range("L111").Value = 2100
$ws.Range("M111").Value = 1153.6666
$ws.Range("N111").Value = -8234

$ws.Range("H118").Value = 801.5714
$ws.Range("I118").Value = 786.46155
$ws.Range("J118").Value = 998
$ws.Range("K118").Value = 2359.38465
$ws.Range("L118").Value = 2994
$ws.Range("M118").Value = -702.38465
$ws.Range("N118").Value = -6308

$ws.Range("H122").Value = 1726.5555
$ws.Range("I122").Value = 1409.6923
$ws.Range("J122").Value = 2550.4
$ws.Range("K122").Value = 4229.0769
$ws.Range("L122").Value = 7651.200000000001
$ws.Range("M122").Value = -1779.0769
$ws.Range("N122").Value = -12551.2

$ws.Range("H132").Value = 2677
$ws.Range("I132").Value = 3007.3333
$ws.Range("J132").Value = 199.5
$ws.Range("K132").Value = 9021.999899999999
$ws.Range("L132").Value = 598.5
$ws.Range("M132").Value = -6491.999899999999
$ws.Range("N132").Value = -5658.5

$ws.Range("H137").Value = 1566.4166
$ws.Range("I137").Value = 1566.4166
$ws.Range("K137").Value = 4699.2498
$ws.Range("M137").Value = -2149.2498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5409.037
$ws.Range("I32").Value = 4630.56
$ws.Range("J32").Value = 15140
$ws.Range("K32").Value = 4630.56
$ws.Range("L32").Value = 15140
$ws.Range("M32").Value = -4343.56
$ws.Range("N32").Value = -15714

$ws.Range("H97").Value = 292.0909
$ws.Range("I97").Value = 223
$ws.Range("J97").Value = 476.33334
$ws.Range("K97").Value = 223
$ws.Range("L97").Value = 476.33334
$ws.Range("M97").Value = 273
$ws.Range("N97").Value = -1468.33334

$ws.Range("H102").Value = 2686.818
$ws.Range("I102").Value = 2455.5
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2455.5
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -833.5
$ws.Range("N102").Value = -8244

$ws.Range("H114").Value = 69666.664
$ws.Range("J114").Value = 69666.664
$ws.Range("L114").Value = 69666.664
$ws.Range("N114").Value = -78344.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2691.8
$ws.Range("I20").Value = 999
$ws.Range("J20").Value = 3115
$ws.Range("K20").Value = 999
$ws.Range("L20").Value = 3115
$ws.Range("M20").Value = -752
$ws.Range("N20").Value = -3609

$ws.Range("H105").Value = 926.3333
$ws.Range("I105").Value = 789.5
$ws.Range("K105").Value = 789.5
$ws.Range("M105").Value = 957.5

$ws.Range("H135").Value = 41999.43
$ws.Range("J135").Value = 41999.43
$ws.Range("L135").Value = 41999.43
$ws.Range("N135").Value = -52139.43

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1086.1
$ws.Range("I11").Value = 1142.1666
$ws.Range("J11").Value = 1002
$ws.Range("K11").Value = 3426.4998
$ws.Range("L11").Value = 3006
$ws.Range("M11").Value = -3286.4998
$ws.Range("N11").Value = -3286

$ws.Range("H14").Value = 249.33333
$ws.Range("I14").Value = 249.33333
$ws.Range("K14").Value = 747.99999
$ws.Range("M14").Value = -574.99999

$ws.Range("H32").Value = 2750
$ws.Range("I32").Value = 550
$ws.Range("J32").Value = 4950
$ws.Range("K32").Value = 1650
$ws.Range("L32").Value = 14850
$ws.Range("M32").Value = -1367
$ws.Range("N32").Value = -15416

$ws.Range("H40").Value = 281.25
$ws.Range("J40").Value = 280
$ws.Range("L40").Value = 1120
$ws.Range("N40").Value = -1258

$ws.Range("H44").Value = 799.5
$ws.Range("I44").Value = 99
$ws.Range("J44").Value = 1500
$ws.Range("K44").Value = 297
$ws.Range("L44").Value = 4500
$ws.Range("M44").Value = 101
$ws.Range("N44").Value = -5296

$ws.Range("H46").Value = 2747
$ws.Range("J46").Value = 4499
$ws.Range("L46").Value = 13497
$ws.Range("N46").Value = -13679

$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -2750
$ws.Range("N48").ClearContents()

$ws.Range("H49").Value = 100
$ws.Range("I49").Value = 100
$ws.Range("K49").Value = 300
$ws.Range("M49").Value = -144

$ws.Range("H57").Value = 6651.6665
$ws.Range("I57").Value = 4977.5
$ws.Range("J57").Value = 10000
$ws.Range("K57").Value = 14932.5
$ws.Range("L57").Value = 30000
$ws.Range("M57").Value = -14373.5
$ws.Range("N57").Value = -31118

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H113").Value = 580.5
$ws.Range("I113").Value = 418.75
$ws.Range("J113").Value = 796.1667
$ws.Range("K113").Value = 1256.25
$ws.Range("L113").Value = 2388.5001
$ws.Range("M113").Value = 913.75
$ws.Range("N113").Value = -6728.5001

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H122").Value = 781.8333
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 474.15
$ws.Range("I2").Value = 497.16666
$ws.Range("K2").Value = 497.16666
$ws.Range("M2").Value = -384.16666

$ws.Range("H43").Value = 25539.572
$ws.Range("J43").Value = 32155.4
$ws.Range("L43").Value = 32155.4
$ws.Range("N43").Value = -32457.4

$ws.Range("H80").Value = 2639.3076
$ws.Range("J80").Value = 3286.4443
$ws.Range("L80").Value = 3286.4443
$ws.Range("N80").Value = -5282.4443

$ws.Range("H83").Value = 2639.3076
$ws.Range("J83").Value = 3286.4443
$ws.Range("L83").Value = 16432.2215
$ws.Range("N83").Value = -26416.2215

$ws.Range("H97").Value = 2395.3333
$ws.Range("I97").Value = 2395.3333
$ws.Range("K97").Value = 2395.3333
$ws.Range("M97").Value = -1899.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 450
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -227
$ws.Range("N55").Value = -846

$ws.Range("H61").Value = 4082.3333
$ws.Range("I61").Value = 3232.1333
$ws.Range("J61").Value = 8333.333000000001
$ws.Range("K61").Value = 3232.1333
$ws.Range("L61").Value = 8333.333000000001
$ws.Range("M61").Value = -3030.1333
$ws.Range("N61").Value = -8737.333000000001

$ws.Range("H82").Value = 1125.1428
$ws.Range("I82").Value = 493.5
$ws.Range("K82").Value = 493.5
$ws.Range("M82").Value = -132.5

$ws.Range("H85").Value = 1125.1428
$ws.Range("I85").Value = 493.5
$ws.Range("K85").Value = 493.5
$ws.Range("M85").Value = 754.5

$ws.Range("H113").Value = 4082.3333
$ws.Range("I113").Value = 3232.1333
$ws.Range("J113").Value = 8333.333000000001
$ws.Range("K113").Value = 3232.1333
$ws.Range("L113").Value = 8333.333000000001
$ws.Range("M113").Value = -1062.1333
$ws.Range("N113").Value = -12673.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 922.73334
$ws.Range("I122").Value = 922.73334
$ws.Range("K122").Value = 2768.20002
$ws.Range("M122").Value = -318.2000200000002

$ws.Range("H132").Value = 1089
$ws.Range("I132").Value = 1089
$ws.Range("K132").Value = 3267
$ws.Range("M132").Value = -737
